# Update countries & provincias Spain
# - Swap display order of "Finlandia" / "Republica de Macedonia" (rows 82/83)
# - Swap display order of "Groenlandia" / "Islas Malvinas" (rows 209/210)
# - Refresh the "Datos actualizados ..." timestamp (row 1)
# - Refresh case statistics for a handful of countries (new scrape snapshot)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "last updated" timestamp (row 1, col A) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Julio de 2020 a las 15:19"

# --- 2. Swap country-name labels (data/stats stay attached to the row, only ---
#        the displayed name moves) ---
$a82 = $ws.Range("A82").Value()
$a83 = $ws.Range("A83").Value()
$ws.Range("A82").Value = $a83
$ws.Range("A83").Value = $a82

$a209 = $ws.Range("A209").Value()
$a210 = $ws.Range("A210").Value()
$ws.Range("A209").Value = $a210
$ws.Range("A210").Value = $a209

# --- 3. Refresh numeric statistics (columns B:H) for the updated rows ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 3097602
$ws.Range("C4").Value = 518
$ws.Range("D4").Value = 1355623
$ws.Range("E4").Value = 1607985
$ws.Range("G4").Value = 22
$ws.Range("H4").Value = 133994

# Row 6 - India
$ws.Range("B6").Value = 746824
$ws.Range("C6").Value = 3343
$ws.Range("D6").Value = 459294
$ws.Range("E6").Value = 266845
$ws.Range("G6").Value = 32
$ws.Range("H6").Value = 20685

# Row 16 - Arabia Saudita
$ws.Range("B16").Value = 220144
$ws.Range("C16").Value = 3036
$ws.Range("D16").Value = 158050
$ws.Range("E16").Value = 60035
$ws.Range("G16").Value = 42
$ws.Range("H16").Value = 2059

# Row 24 - Catar
$ws.Range("B24").Value = 101553
$ws.Range("C24").Value = 608
$ws.Range("D24").Value = 96107
$ws.Range("E24").Value = 5308
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 138

# Row 26 - Argentina
$ws.Range("D26").Value = 36502
$ws.Range("E26").Value = 45270
$ws.Range("G26").Value = 10
$ws.Range("H26").Value = 1654

# Row 30 - Irak
$ws.Range("B30").Value = 67442
$ws.Range("C30").Value = 2741
$ws.Range("D30").Value = 37879
$ws.Range("E30").Value = 26784
$ws.Range("G30").Value = 94
$ws.Range("H30").Value = 2779

# Row 31 - Bielorrusia
$ws.Range("B31").Value = 64224
$ws.Range("C31").Value = 221
$ws.Range("D31").Value = 52854
$ws.Range("E31").Value = 10927
$ws.Range("G31").Value = 7
$ws.Range("H31").Value = 443

# Row 37 - Paises Bajos
$ws.Range("B37").Value = 50746
$ws.Range("C37").Value = 52
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 6135

# Row 68 - Dinamarca
$ws.Range("B68").Value = 12900
$ws.Range("C68").Value = 12
$ws.Range("D68").Value = 12001
$ws.Range("E68").Value = 290

# Row 71 - Uzbekistan
$ws.Range("E71").Value = 3985
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 42

# Row 82 (label now "Republica de Macedonia" after the swap above)
$ws.Range("B82").Value = 7406
$ws.Range("C82").Value = 162
$ws.Range("D82").Value = 3554
$ws.Range("E82").Value = 3493
$ws.Range("G82").Value = 8
$ws.Range("H82").Value = 359

# Row 83 (label now "Finlandia" after the swap above)
$ws.Range("B83").Value = 7265
$ws.Range("C83").Value = 3
$ws.Range("D83").Value = 6800
$ws.Range("E83").Value = 136
$ws.Range("H83").Value = 329

# Row 101 - Croacia
$ws.Range("B101").Value = 3325
$ws.Range("C101").Value = 53
$ws.Range("D101").Value = 2277
$ws.Range("E101").Value = 934
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 114

# Row 110 - Cuba
$ws.Range("B110").Value = 2399
$ws.Range("C110").Value = 4
$ws.Range("D110").Value = 2242
$ws.Range("E110").Value = 71
